$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.760.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.549.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.96%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.21%  "
$ws.Range("E11").Value = "  -4.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.114"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.939.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.543.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.01%  "
$ws.Range("E16").Value = "  -6.01%  "
$ws.Range("E17").Value = "  -6.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.763.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("E25").Value = "  -5.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.29%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.25%  "
$ws.Range("E35").Value = "  -10.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0794"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.80%  "
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.20%  "
$ws.Range("E41").Value = "  -4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.084.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -10.62%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.17%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.792.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("E51").Value = "  -4.00%  "

Write-Host "Applied 94 cell updates"
